$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 79 : 2022-05-25, 08:50 -> 09:30, "Scéance Tests" description + training testers
$ws.Range("A79").Value = [DateTime]"2022-05-25"
$ws.Range("B79").Value = [DateTime]"2022-05-25 08:50:00"
$ws.Range("C79").Value = [DateTime]"2022-05-25 09:30:00"
$ws.Range("E79").Value = "Documentation: Stratégie de test + Réalisation "

# Row 80 : 2022-05-25, 09:50 -> 10:00
$ws.Range("A80").Value = [DateTime]"2022-05-25"
$ws.Range("B80").Value = [DateTime]"2022-05-25 09:50:00"
$ws.Range("C80").Value = [DateTime]"2022-05-25 10:00:00"
$ws.Range("E80").Value = "Scéance Tests"
$ws.Range("F80").Value = "Testeurs:`nAzad Saffai`nThirusan Rajadurai"

# Row 81 : 2022-05-25, 10:00 -> 12:15
$ws.Range("A81").Value = [DateTime]"2022-05-25"
$ws.Range("B81").Value = [DateTime]"2022-05-25 10:00:00"
$ws.Range("C81").Value = [DateTime]"2022-05-25 12:15:00"
$ws.Range("E81").Value = "Documentation Fin réalisation + tests effectué"

# Row 82 : 2022-05-25, 13:30 -> 16:55
$ws.Range("A82").Value = [DateTime]"2022-05-25"
$ws.Range("B82").Value = [DateTime]"2022-05-25 13:30:00"
$ws.Range("C82").Value = [DateTime]"2022-05-25 16:55:00"
$ws.Range("E82").Value = "Documentation: Conclusion + résumé"

# Update visible view window / active selection to match final editor state
$ws.Application.ActiveWindow.ScrollRow = 79
$ws.Range("E83").Select()

$wb.Application.CalculateFull()
